# edit.ps1
#
# Applies three related changes to the document:
#
#   1. The title paragraph's paragraph-mark run properties (w:pPr/w:rPr) lose
#      their <w:rFonts w:hint="cs"/> element (the run itself keeps it).
#   2. The text "תמחקי את זה אחר כך" is shortened to "תמחקי את זה "
#      (trailing "אחר כך" removed, trailing space kept).
#   3. The hidden "_GoBack" bookmark moves from the following (empty)
#      paragraph to right after the shortened run, inside the same
#      paragraph as the shortened text.

$d = $word.ActiveDocument

# --- 1. Strip the rFonts hint from the first paragraph's paragraph mark ----
# The paragraph mark's own rPr (inside pPr) should no longer carry
# <w:rFonts w:hint="cs"/>, while the run's rPr is left untouched.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00385FFE" w:rsidRDefault="00385FFE" w:rsidP="00385FFE">
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
                <w:u w:val="single"/>
                <w:rtl/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00385FFE">
              <w:rPr>
                <w:rFonts w:hint="cs"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
                <w:u w:val="single"/>
                <w:rtl/>
              </w:rPr>
              <w:t>סיפור לקוח</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$titleRange.InsertXML($titleXml)

# --- 2 & 3. Move the _GoBack bookmark and shorten the reminder text -------
# Locate the paragraph containing "תמחקי את זה אחר כך" by searching for it,
# then compute where the new bookmark / cut point belongs: right after the
# text that is being kept ("תמחקי את זה ", 12 characters incl. the space).
$reminderRange = $d.Content
$reminderRange.Find.ClearFormatting()
$found = $reminderRange.Find.Execute("תמחקי את זה אחר כך", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $keepText = "תמחקי את זה "
    $cutPoint = $reminderRange.Start + $keepText.Length

    # Remove the old (hidden) _GoBack bookmark first so the name is free.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    # Re-create _GoBack at the point where the kept text ends, while the
    # paragraph is still long -- this must happen before the trailing text
    # is deleted, while the insertion point sits safely inside the run.
    $bookmarkRange = $d.Range($cutPoint, $cutPoint)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    # Now delete the trailing "אחר כך" that followed the kept text.
    $tailRange = $d.Range($cutPoint, $reminderRange.End)
    $tailRange.Text = ""
}
